$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 126038.19
$ws.Range("I28").Value = 143793.64
$ws.Range("J28").Value = 1750
$ws.Range("K28").Value = 143793.64
$ws.Range("L28").Value = 1750
$ws.Range("M28").Value = -143308.64
$ws.Range("N28").Value = -2720
$ws.Range("H39").Value = 635.1667
$ws.Range("I39").Value = 635.1667
$ws.Range("K39").Value = 1905.5001
$ws.Range("M39").Value = -1609.5001
$ws.Range("H53").Value = 1292.3846
$ws.Range("I53").Value = 757.1667
$ws.Range("K53").Value = 757.1667
$ws.Range("M53").Value = -120.1667
$ws.Range("H74").Value = 4952.875
$ws.Range("I74").Value = 4917.5713
$ws.Range("K74").Value = 4917.5713
$ws.Range("M74").Value = -3981.5713
$ws.Range("H77").Value = 4952.875
$ws.Range("I77").Value = 4917.5713
$ws.Range("K77").Value = 24587.8565
$ws.Range("M77").Value = -19907.8565
$ws.Range("H88").Value = 3624.4443
$ws.Range("J88").Value = 3487.3333
$ws.Range("L88").Value = 3487.3333
$ws.Range("N88").Value = -4299.3333
$ws.Range("H91").Value = 3624.4443
$ws.Range("J91").Value = 3487.3333
$ws.Range("L91").Value = 3487.3333
$ws.Range("N91").Value = -6295.3333
$ws.Range("H100").Value = 4169.357
$ws.Range("I100").Value = 3249.8333
$ws.Range("K100").Value = 3249.8333
$ws.Range("M100").Value = -2708.8333
$ws.Range("H106").Value = 5130475
$ws.Range("I106").Value = 5130475
$ws.Range("K106").Value = 5130475
$ws.Range("M106").Value = -5129844
$ws.Range("H138").Value = 2310.1
$ws.Range("I138").Value = 958.875
$ws.Range("J138").Value = 2736.8027
$ws.Range("K138").Value = 2876.625
$ws.Range("L138").Value = 8210.408100000001
$ws.Range("M138").Value = 2263.375
$ws.Range("N138").Value = -18490.4081

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27444368
$ws.Range("I32").Value = 29612026
$ws.Range("J32").Value = 11909480
$ws.Range("K32").Value = 29612026
$ws.Range("L32").Value = 11909480
$ws.Range("M32").Value = -29611739
$ws.Range("N32").Value = -11910054
$ws.Range("H45").Value = 2936.0833
$ws.Range("I45").Value = 1233.25
$ws.Range("K45").Value = 1233.25
$ws.Range("M45").Value = -856.25
$ws.Range("H74").Value = 2895.2903
$ws.Range("I74").Value = 2353.889
$ws.Range("K74").Value = 2353.889
$ws.Range("M74").Value = -1479.889
$ws.Range("H77").Value = 2895.2903
$ws.Range("I77").Value = 2353.889
$ws.Range("K77").Value = 11769.445
$ws.Range("M77").Value = -7401.445
$ws.Range("H88").Value = 1254.5555
$ws.Range("J88").Value = 1156
$ws.Range("L88").Value = 1156
$ws.Range("N88").Value = -1968
$ws.Range("H91").Value = 1254.5555
$ws.Range("J91").Value = 1156
$ws.Range("L91").Value = 1156
$ws.Range("N91").Value = -3964
$ws.Range("H122").Value = 3734.2354
$ws.Range("I122").Value = 3343.6667
$ws.Range("K122").Value = 10031.0001
$ws.Range("M122").Value = -7581.000100000001
$ws.Range("H128").Value = 62244.25
$ws.Range("J128").Value = 62244.25
$ws.Range("L128").Value = 62244.25
$ws.Range("N128").Value = -72204.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2549.25
$ws.Range("I86").Value = 2280.4
$ws.Range("K86").Value = 2280.4
$ws.Range("M86").Value = -1157.4
$ws.Range("H89").Value = 2549.25
$ws.Range("I89").Value = 2280.4
$ws.Range("K89").Value = 11402
$ws.Range("M89").Value = -5786
$ws.Range("H105").Value = 2749.8262
$ws.Range("I105").Value = 2364.5
$ws.Range("K105").Value = 2364.5
$ws.Range("M105").Value = -617.5
$ws.Range("H107").Value = 1582.6666
$ws.Range("J107").Value = 1632.3334
$ws.Range("L107").Value = 1632.3334
$ws.Range("N107").Value = -5472.3334
$ws.Range("H134").Value = 1882988.8
$ws.Range("I134").Value = 2465151
$ws.Range("K134").Value = 7395453
$ws.Range("M134").Value = -7392918

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1565.8334
$ws.Range("I16").Value = 1273.75
$ws.Range("J16").Value = 2150
$ws.Range("K16").Value = 1273.75
$ws.Range("L16").Value = 2150
$ws.Range("M16").Value = -986.75
$ws.Range("N16").Value = -2724
$ws.Range("H31").Value = 5603.5186
$ws.Range("I31").Value = 2268.4375
$ws.Range("J31").Value = 10454.546
$ws.Range("K31").Value = 2268.4375
$ws.Range("L31").Value = 10454.546
$ws.Range("M31").Value = -1973.4375
$ws.Range("N31").Value = -11044.546
$ws.Range("H34").Value = 5603.5186
$ws.Range("I34").Value = 2268.4375
$ws.Range("J34").Value = 10454.546
$ws.Range("K34").Value = 2268.4375
$ws.Range("L34").Value = 10454.546
$ws.Range("M34").Value = -2066.4375
$ws.Range("N34").Value = -10858.546
$ws.Range("H58").Value = 2752.532
$ws.Range("I58").Value = 2514.8
$ws.Range("J58").Value = 4111
$ws.Range("K58").Value = 2514.8
$ws.Range("L58").Value = 4111
$ws.Range("M58").Value = -2311.8
$ws.Range("N58").Value = -4517
$ws.Range("H99").Value = 2467.75
$ws.Range("I99").Value = 2290.3333
$ws.Range("K99").Value = 2290.3333
$ws.Range("M99").Value = -792.3332999999998
$ws.Range("H105").Value = 1928.091
$ws.Range("I105").Value = 1291.8
$ws.Range("K105").Value = 1291.8
$ws.Range("M105").Value = 455.2
$ws.Range("H107").Value = 1183.4445
$ws.Range("I107").Value = 543.9286
$ws.Range("K107").Value = 543.9286
$ws.Range("M107").Value = 1376.0714
$ws.Range("H113").Value = 1565.8334
$ws.Range("I113").Value = 1273.75
$ws.Range("J113").Value = 2150
$ws.Range("K113").Value = 1273.75
$ws.Range("L113").Value = 2150
$ws.Range("M113").Value = 896.25
$ws.Range("N113").Value = -6490
$ws.Range("H125").Value = 99974.5
$ws.Range("J125").Value = 99974.5
$ws.Range("L125").Value = 99974.5
$ws.Range("N125").Value = -104894.5
$ws.Range("H126").Value = 2467.75
$ws.Range("I126").Value = 2290.3333
$ws.Range("K126").Value = 6870.999899999999
$ws.Range("M126").Value = -4400.999899999999
$ws.Range("H136").Value = 2752.532
$ws.Range("I136").Value = 2514.8
$ws.Range("J136").Value = 4111
$ws.Range("K136").Value = 7544.400000000001
$ws.Range("L136").Value = 12333
$ws.Range("M136").Value = -4994.400000000001
$ws.Range("N136").Value = -17433

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3374.9285
$ws.Range("I109").Value = 2178.4285
$ws.Range("J109").Value = 4571.4287
$ws.Range("K109").Value = 6535.2855
$ws.Range("L109").Value = 13714.2861
$ws.Range("M109").Value = -5495.2855
$ws.Range("N109").Value = -15794.2861
$ws.Range("H137").Value = 3732.0557
$ws.Range("I137").Value = 1191.6666
$ws.Range("J137").Value = 5002.25
$ws.Range("K137").Value = 3574.9998
$ws.Range("L137").Value = 15006.75
$ws.Range("M137").Value = 1525.0002
$ws.Range("N137").Value = -25206.75
$ws.Range("H138").Value = 21301520
$ws.Range("I138").Value = 1899.7142
$ws.Range("K138").Value = 5699.142599999999
$ws.Range("M138").Value = -559.1425999999992

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 102989.5
$ws.Range("J117").Value = 102989.5
$ws.Range("L117").Value = 102989.5
$ws.Range("N117").Value = -109873.5
$ws.Range("H122").Value = 2183
$ws.Range("J122").Value = 2183
$ws.Range("L122").Value = 6549
$ws.Range("N122").Value = -11449
$ws.Range("H132").Value = 4226.727
$ws.Range("I132").Value = 3311.875
$ws.Range("K132").Value = 9935.625
$ws.Range("M132").Value = -7405.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25645288
$ws.Range("I40").Value = 41669590
$ws.Range("K40").Value = 41669590
$ws.Range("M40").Value = -41669454
$ws.Range("H100").Value = 2810.5
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("H115").Value = 73495
$ws.Range("J115").Value = 73495
$ws.Range("L115").Value = 73495
$ws.Range("N115").Value = -75845
$ws.Range("H122").Value = 36222.668
$ws.Range("I122").Value = 37625.5
$ws.Range("J122").Value = 25000
$ws.Range("K122").Value = 112876.5
$ws.Range("L122").Value = 75000
$ws.Range("M122").Value = -110426.5
$ws.Range("N122").Value = -79900
$ws.Range("H131").Value = 107996.336
$ws.Range("J131").Value = 107996.336
$ws.Range("L131").Value = 107996.336
$ws.Range("N131").Value = -118076.336
$ws.Range("H136").Value = 2271.2856
$ws.Range("I136").Value = 2079.8
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 6239.400000000001
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -3689.400000000001
$ws.Range("N136").Value = -13350
$ws.Range("N100").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3661.889
$ws.Range("I81").Value = 3101.5557
$ws.Range("J81").Value = 4222.222
$ws.Range("K81").Value = 6203.1114
$ws.Range("L81").Value = 8444.444
$ws.Range("M81").Value = -5142.1114
$ws.Range("N81").Value = -10566.444
$ws.Range("H84").Value = 3661.889
$ws.Range("I84").Value = 3101.5557
$ws.Range("J84").Value = 4222.222
$ws.Range("K84").Value = 31015.557
$ws.Range("L84").Value = 42222.22
$ws.Range("M84").Value = -25711.557
$ws.Range("N84").Value = -52830.22
$ws.Range("H120").Value = 87321.86
$ws.Range("J120").Value = 87321.86
$ws.Range("L120").Value = 87321.86
$ws.Range("N120").Value = -96997.86
$ws.Range("H122").Value = 166675000
$ws.Range("I122").Value = 166675000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 500025000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -500022550
$ws.Range("N122").ClearContents()
